$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.202.47'
$ws.Range('E2').Value = '  -3.45%  '
$ws.Range('D3').Value = '2.465.19'
$ws.Range('E3').Value = '  -2.39%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'311.95"
$ws.Range('E5').Value = '  +1.03%  '
$ws.Range('D6').Value = "'94.16"
$ws.Range('E6').Value = '  -5.91%  '
$ws.Range('D7').Value = "'0.552"
$ws.Range('E7').Value = '  -2.71%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = "'0.498"
$ws.Range('E9').Value = '  -4.24%  '
$ws.Range('D10').Value = "'33.33"
$ws.Range('E10').Value = '  -6.02%  '
$ws.Range('D11').Value = "'0.0778"
$ws.Range('E11').Value = '  -3.22%  '
$ws.Range('E12').Value = '  -1.09%  '
$ws.Range('D13').Value = "'6.98"
$ws.Range('E13').Value = '  -4.22%  '
$ws.Range('D14').Value = '2.842.64'
$ws.Range('E14').Value = '  -2.56%  '
$ws.Range('D15').Value = '2.455.41'
$ws.Range('E15').Value = '  -6.80%  '
$ws.Range('D16').Value = "'14.85"
$ws.Range('E16').Value = '  -2.75%  '
$ws.Range('D17').Value = "'0.785"
$ws.Range('E17').Value = '  -3.05%  '
$ws.Range('D18').Value = '41.151.60'
$ws.Range('E18').Value = '  -3.51%  '
$ws.Range('D19').Value = "'6.30"
$ws.Range('E19').Value = '  -5.53%  '
$ws.Range('D20').Value = '0.0₃0921'
$ws.Range('E20').Value = '  -2.74%  '
$ws.Range('D21').Value = "'11.30"
$ws.Range('E21').Value = '  -7.32%  '
$ws.Range('D22').Value = "'68.47"
$ws.Range('E22').Value = '  -1.04%  '
$ws.Range('D23').Value = "'235.70"
$ws.Range('E23').Value = '  -2.65%  '
$ws.Range('D24').Value = "'2.75"
$ws.Range('E24').Value = '  -3.36%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').Value = "'1.91"
$ws.Range('E26').Value = '  -5.55%  '
$ws.Range('D27').Value = "'23.99"
$ws.Range('E27').Value = '  -5.40%  '
$ws.Range('E28').Value = '  -6.17%  '
$ws.Range('D29').Value = "'9.61"
$ws.Range('E29').Value = '  -5.01%  '
$ws.Range('D30').Value = "'36.51"
$ws.Range('E30').Value = '  -4.69%  '
$ws.Range('D31').Value = "'152.72"
$ws.Range('E31').Value = '  -4.27%  '
$ws.Range('D32').Value = "'5.49"
$ws.Range('E32').Value = '  -4.29%  '
$ws.Range('E33').Value = '  -5.37%  '
$ws.Range('E34').Value = '  -3.36%  '
$ws.Range('D35').Value = "'0.0743"
$ws.Range('E35').Value = '  -4.62%  '
$ws.Range('D36').Value = "'3.05"
$ws.Range('E36').Value = '  -1.84%  '
$ws.Range('D37').Value = "'17.10"
$ws.Range('E37').Value = '  -7.04%  '
$ws.Range('D38').Value = "'1.88"
$ws.Range('E38').Value = '  -3.80%  '
$ws.Range('E39').Value = '  -2.83%  '
$ws.Range('D40').Value = "'0.102"
$ws.Range('E40').Value = '  -7.39%  '
$ws.Range('D41').Value = "'4.17"
$ws.Range('E41').Value = '  -1.09%  '
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('D43').Value = "'20.15"
$ws.Range('E43').Value = '  -9.30%  '
$ws.Range('D44').Value = '1.969.54'
$ws.Range('E44').Value = '  -1.48%  '
$ws.Range('D45').Value = "'0.0284"
$ws.Range('E45').Value = '  -5.00%  '
$ws.Range('E46').Value = '  -7.41%  '
$ws.Range('E47').Value = '  -1.78%  '
$ws.Range('B48').Value = 'ordi'
$ws.Range('C48').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D48').Value = "'69.27"
$ws.Range('E48').Value = '  -3.04%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = "'97.10"
$ws.Range('E49').Value = '  -3.46%  '
$ws.Range('E50').Value = '  -5.86%  '
$ws.Range('D51').Value = "'73.80"
$ws.Range('E51').Value = '  -6.85%  '
